$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1261.6428
$ws.Range("I17").Value = 735
$ws.Range("J17").Value = 1302.1538
$ws.Range("K17").Value = 2205
$ws.Range("L17").Value = 3906.4614
$ws.Range("M17").Value = -2037
$ws.Range("N17").Value = -4242.4614
$ws.Range("H87").Value = 48916.332
$ws.Range("J87").Value = 48916.332
$ws.Range("L87").Value = 48916.332
$ws.Range("N87").Value = -51412.332
$ws.Range("H90").Value = 48916.332
$ws.Range("J90").Value = 48916.332
$ws.Range("L90").Value = 146748.996
$ws.Range("N90").Value = -159228.996
$ws.Range("H98").Value = 37039960
$ws.Range("I98").Value = 38464384
$ws.Range("K98").Value = 38464384
$ws.Range("M98").Value = -38462886
$ws.Range("H122").Value = 37039960
$ws.Range("I122").Value = 38464384
$ws.Range("K122").Value = 115393152
$ws.Range("M122").Value = -115390702
$ws.Range("H132").Value = 1643.5
$ws.Range("I132").Value = 1643.5
$ws.Range("K132").Value = 4930.5
$ws.Range("M132").Value = -2400.5
$ws.Range("H138").Value = 2277645
$ws.Range("J138").Value = 3131088.2
$ws.Range("L138").Value = 9393264.600000001
$ws.Range("N138").Value = -9403544.600000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 21307.375
$ws.Range("I28").Value = 2563.4
$ws.Range("K28").Value = 2563.4
$ws.Range("M28").Value = -2371.4
$ws.Range("H43").Value = 28395.334
$ws.Range("J43").Value = 28395.334
$ws.Range("L43").Value = 28395.334
$ws.Range("N43").Value = -29021.334
$ws.Range("H61").Value = 32263534
$ws.Range("I61").Value = 1561
$ws.Range("K61").Value = 1561
$ws.Range("M61").Value = -1349
$ws.Range("H88").Value = 1420
$ws.Range("I88").Value = 920
$ws.Range("J88").Value = 1720
$ws.Range("K88").Value = 920
$ws.Range("L88").Value = 1720
$ws.Range("M88").Value = -514
$ws.Range("N88").Value = -2532
$ws.Range("H91").Value = 1420
$ws.Range("I91").Value = 920
$ws.Range("J91").Value = 1720
$ws.Range("K91").Value = 920
$ws.Range("L91").Value = 1720
$ws.Range("M91").Value = 484
$ws.Range("N91").Value = -4528
$ws.Range("H97").Value = 6411979
$ws.Range("I97").Value = 1431
$ws.Range("J97").Value = 20835712
$ws.Range("K97").Value = 1431
$ws.Range("L97").Value = 20835712
$ws.Range("M97").Value = -935
$ws.Range("N97").Value = -20836704
$ws.Range("H99").Value = 21307.375
$ws.Range("I99").Value = 2563.4
$ws.Range("K99").Value = 2563.4
$ws.Range("M99").Value = 431.5999999999999
$ws.Range("H122").Value = 4385.4546
$ws.Range("I122").Value = 3289.158
$ws.Range("K122").Value = 9867.474
$ws.Range("M122").Value = -7417.474
$ws.Range("H132").Value = 5202.25
$ws.Range("I132").Value = 2982.6487
$ws.Range("K132").Value = 8947.946100000001
$ws.Range("M132").Value = -6417.946100000001
$ws.Range("H136").Value = 32263534
$ws.Range("I136").Value = 1561
$ws.Range("K136").Value = 4683
$ws.Range("M136").Value = -2133

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5480.603
$ws.Range("I31").Value = 1986.6552
$ws.Range("K31").Value = 1986.6552
$ws.Range("M31").Value = -1691.6552
$ws.Range("H34").Value = 5480.603
$ws.Range("I34").Value = 1986.6552
$ws.Range("K34").Value = 1986.6552
$ws.Range("M34").Value = -1784.6552
$ws.Range("H132").Value = 8406.5625
$ws.Range("I132").Value = 4302
$ws.Range("J132").Value = 10272.272
$ws.Range("K132").Value = 12906
$ws.Range("L132").Value = 30816.816
$ws.Range("M132").Value = -10376
$ws.Range("N132").Value = -35876.81600000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6620705
$ws.Range("I4").Value = 10105134
$ws.Range("K4").Value = 30315402
$ws.Range("M4").Value = -30315290
$ws.Range("H114").Value = 671.0714
$ws.Range("I114").Value = 678.6
$ws.Range("K114").Value = 2035.8
$ws.Range("M114").Value = 1218.2
$ws.Range("H117").Value = 55556736
$ws.Range("I117").Value = 23810550
$ws.Range("J117").Value = 71429830
$ws.Range("K117").Value = 71431650
$ws.Range("L117").Value = 214289490
$ws.Range("M117").Value = -71428208
$ws.Range("N117").Value = -214296374
$ws.Range("H121").Value = 3847639.2
$ws.Range("I121").Value = 1491.4
$ws.Range("J121").Value = 16668132
$ws.Range("K121").Value = 4474.200000000001
$ws.Range("L121").Value = 50004396
$ws.Range("M121").Value = -3164.200000000001
$ws.Range("N121").Value = -50007016
$ws.Range("H129").Value = 35786572
$ws.Range("I129").Value = 834
$ws.Range("J129").Value = 100200904
$ws.Range("K129").Value = 2502
$ws.Range("L129").Value = 300602712
$ws.Range("M129").Value = 2498
$ws.Range("N129").Value = -300612712

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 92.77778000000001
$ws.Range("I2").Value = 112.166664
$ws.Range("J2").Value = 54
$ws.Range("K2").Value = 112.166664
$ws.Range("L2").Value = 54
$ws.Range("M2").Value = 0.8333360000000027
$ws.Range("N2").Value = -280
$ws.Range("H15").Value = 7675.3335
$ws.Range("J15").Value = 7675.3335
$ws.Range("L15").Value = 7675.3335
$ws.Range("N15").Value = -8251.333500000001
$ws.Range("H81").Value = 7675.3335
$ws.Range("J81").Value = 7675.3335
$ws.Range("L81").Value = 7675.3335
$ws.Range("N81").Value = -9671.333500000001
$ws.Range("H84").Value = 7675.3335
$ws.Range("J84").Value = 7675.3335
$ws.Range("L84").Value = 23026.0005
$ws.Range("N84").Value = -33010.00049999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4711.2856
$ws.Range("I7").Value = 3024.8235
$ws.Range("J7").Value = 7317.636
$ws.Range("K7").Value = 3024.8235
$ws.Range("L7").Value = 7317.636
$ws.Range("M7").Value = -2912.8235
$ws.Range("N7").Value = -7541.636
$ws.Range("H93").Value = 5063.5
$ws.Range("I93").Value = 5097.25
$ws.Range("J93").Value = 4996
$ws.Range("K93").Value = 5097.25
$ws.Range("L93").Value = 4996
$ws.Range("M93").Value = -3849.25
$ws.Range("N93").Value = -7492
$ws.Range("H100").Value = 5183.8887
$ws.Range("I100").Value = 2641.5
$ws.Range("K100").Value = 2641.5
$ws.Range("M100").Value = -2100.5
$ws.Range("H126").Value = 4711.2856
$ws.Range("I126").Value = 3024.8235
$ws.Range("J126").Value = 7317.636
$ws.Range("K126").Value = 9074.470499999999
$ws.Range("L126").Value = 21952.908
$ws.Range("M126").Value = -6604.470499999999
$ws.Range("N126").Value = -26892.908
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -85060
$ws.Range("H136").Value = 9097.368
$ws.Range("I136").Value = 3499.5
$ws.Range("J136").Value = 11681
$ws.Range("K136").Value = 10498.5
$ws.Range("L136").Value = 35043
$ws.Range("M136").Value = -7948.5
$ws.Range("N136").Value = -40143
